# Weekly driver report update for 2025-04-20
# Applies the data refresh to the "Driver Summary" sheet:
#  - Updates the Bad Drivers Good Roaming Calculation value
#  - Rolls the Good Drivers table: a new driver entry is added at the top
#    (row 12), all the other rows shift down by one, and the oldest
#    entry (previously row 24) drops out of the printed set while new
#    "Driver Vintage" dates are populated for rows 12-23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bad Drivers table -----------------------------------------------
$ws.Range("D3").Value = 98.3

# --- Good Drivers table (rows 12-23 get refreshed) --------------------
# Note: column E ("Driver Vintage") holds date-formatted text, not a real
# date value, so a leading apostrophe is used to force the cell to stay
# text instead of being auto-converted into a date serial number.

$ws.Range("A12").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 23.70.2.3"
$ws.Range("B12").Value = 18721
$ws.Range("D12").Value = 99.90000000000001
$ws.Range("E12").Value = "'2024-07-23"

$ws.Range("A13").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.250.10.1"
$ws.Range("B13").Value = 69578
$ws.Range("D13").Value = 99.90000000000001
$ws.Range("E13").Value = "'2023-08-14"

$ws.Range("A14").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.230.0.8"
$ws.Range("B14").Value = 329845
$ws.Range("D14").Value = 99.90000000000001
$ws.Range("E14").Value = "'2023-05-08"

$ws.Range("A15").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.200.0.6"
$ws.Range("B15").Value = 143808
$ws.Range("D15").Value = 99.90000000000001
$ws.Range("E15").Value = "'2023-01-16"

$ws.Range("A16").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.190.0.4"
$ws.Range("B16").Value = 287148
$ws.Range("D16").Value = 99.90000000000001
$ws.Range("E16").Value = "'2022-11-22"

$ws.Range("A17").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.160.0.4"
$ws.Range("B17").Value = 96526
$ws.Range("D17").Value = 99.90000000000001
$ws.Range("E17").Value = "'2022-08-13"

$ws.Range("A18").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.30.0.11"
$ws.Range("B18").Value = 67111
$ws.Range("D18").Value = 100
$ws.Range("E18").Value = "'2021-01-19"

$ws.Range("A19").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.10.0.7"
$ws.Range("B19").Value = 66577
$ws.Range("D19").Value = 100
$ws.Range("E19").Value = "'2020-10-19"

$ws.Range("A20").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.0.1.1"
$ws.Range("B20").Value = 15734
$ws.Range("D20").Value = 99.90000000000001
$ws.Range("E20").Value = "'2020-09-28"

$ws.Range("A21").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 21.60.2.1"
$ws.Range("B21").Value = 26241
$ws.Range("D21").Value = 100
$ws.Range("E21").Value = "'2019-12-14"

$ws.Range("A22").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 21.40.2.2"
$ws.Range("B22").Value = 88435
$ws.Range("D22").Value = 99.90000000000001
$ws.Range("E22").Value = "'2019-08-31"

$ws.Range("A23").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 21.30.4.1"
$ws.Range("B23").Value = 13016
$ws.Range("D23").Value = 100
$ws.Range("E23").Value = "'2019-07-29"

# Row 24 (Intel(R) Wi-Fi 6 AX200 160MHz - 21.10.1.2) is unchanged.
